# Update catch dependency factors (Table 13) for the new catch-dependency
# inputs. Each row is a country; column 2 = Cod, column 3 = Hake.
# Using Cell(row, col) addressing (rather than a blind text Find/Replace)
# avoids value collisions where a new value equals another cell's old
# value (e.g. 0.971 -> 0.942, while 0.942 -> 0.885 elsewhere).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $newValue) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.MoveEnd(1, -1)  # wdCharacter: drop trailing cell-mark character
    $r.Text = $newValue
}

# Row 2 = BE
Set-CellValue $t 2 2 "0.973"
Set-CellValue $t 2 3 "0.942"

# Row 3 = DK
Set-CellValue $t 3 2 "0.714"
Set-CellValue $t 3 3 "0.8"

# Row 4 = DE
Set-CellValue $t 4 2 "0.947"
Set-CellValue $t 4 3 "0.979"

# Row 5 = EE
Set-CellValue $t 5 2 "0.735"

# Row 6 = IE
Set-CellValue $t 6 2 "0.793"

# Row 7 = ES
Set-CellValue $t 7 2 "0.770"
Set-CellValue $t 7 3 "0.806"

# Row 8 = FR
Set-CellValue $t 8 2 "0.885"
Set-CellValue $t 8 3 "0.903"

# Row 9 = LV
Set-CellValue $t 9 2 "0.679"

# Row 10 = LT
Set-CellValue $t 10 2 "0.693"

# Row 11 = NL
Set-CellValue $t 11 2 "0.642"

# Row 12 = PL
Set-CellValue $t 12 2 "0.760"

# Row 13 = PT
Set-CellValue $t 13 2 "0.461"
Set-CellValue $t 13 3 "0.607"

# Row 14 = FI
Set-CellValue $t 14 2 "0.826"

# Row 15 = SE
Set-CellValue $t 15 2 "0.759"
Set-CellValue $t 15 3 "0.137"
